# Delete row 21 from the "NEW" sheet, shifting all subsequent rows up by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

$ws.Rows.Item(21).Delete()
